$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the quarterly header labels (row 8: fiscal period) ---
# The oldest quarter column (D) label is dropped, everything shifts one
# column to the left, and a brand-new quarter label is appended at M.
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Update the "publish date" row (row 9) the same way ---
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-02-06 (9)"
$ws.Range("F9").Value = "1401-04-29 (3)"
$ws.Range("G9").Value = "1401-08-25 (4)"
$ws.Range("H9").Value = "1401-10-28 (2)"
$ws.Range("I9").Value = "1402-02-10 (8)"
$ws.Range("J9").Value = "1401-04-29"
$ws.Range("K9").Value = "1401-08-25 (2)"
$ws.Range("L9").Value = "1401-10-28"
$ws.Range("M9").Value = "1402-02-10 (2)"

# --- Update the financial data rows: drop oldest quarter (col D), shift
#     everything left, append newest quarter values at column M ---

function Set-RowValues {
    param($ws, $row, $values)
    $cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

# Row 11: فروش (sales)
Set-RowValues $ws 11 @(3584, 4175, 7379, 7879, 8062, 8121, 10680, 8986, 7701, 5347)

# Row 12: بهای تمام شده کالای فروش رفته (COGS)
Set-RowValues $ws 12 @(-2716, -2992, -5427, -5933, -6430, -7012, -8454, -6997, -6191, -4348)

# Row 13: سود (زیان) ناخالص (gross profit)
Set-RowValues $ws 13 @(868, 1183, 1952, 1946, 1632, 1109, 2226, 1989, 1510, 999)

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expense)
Set-RowValues $ws 14 @(-104, -176, -137, -151, -114, -218, -156, -138, -144, -170)

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
Set-RowValues $ws 16 @(26, 95, 109, 125, 108, 97, 65, 53, 135, -108)

# Row 17: سود (زیان) عملیاتی (operating profit)
Set-RowValues $ws 17 @(789, 1102, 1924, 1920, 1626, 989, 2134, 1904, 1501, 722)

# Row 18: هزینه های مالی (financial expense)
Set-RowValues $ws 18 @(-14, -1, -13, -20, -18, -3, -15, 0, -3, -18)

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
Set-RowValues $ws 19 @(13, -225, 1, 65, 15, -15, 7, 23, 29, -71)

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
Set-RowValues $ws 20 @(788, 876, 1911, 1965, 1622, 970, 2126, 1927, 1528, 633)

# Row 21: مالیات (tax)
Set-RowValues $ws 21 @(-105, -6, -430, -259, -403, 1, -415, -214, -316, 194)

# Row 22: سود (زیان) خالص عملیات در حال تداوم
Set-RowValues $ws 22 @(683, 870, 1481, 1706, 1219, 972, 1711, 1713, 1211, 827)

# Row 24: سود (زیان) خالص
Set-RowValues $ws 24 @(683, 870, 1481, 1706, 1219, 972, 1711, 1713, 1211, 827)

# Row 26: سرمایه (capital)
Set-RowValues $ws 26 @(1746, 1943, 2034, 1818, 1667, 1723, 1615, 3063, 2734, 2090)
